$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update simulation-derived transition probabilities in the matrix.
# Values correspond directly to the updated dataset (more simulated games).
# Row 2
$ws.Range("B2").Value = 0.2420382165605096
$ws.Range("C2").Value = 0.4872611464968153
$ws.Range("J2").Value = 0.02547770700636943
$ws.Range("O2").Value = 0.003184713375796179
$ws.Range("P2").Value = 0.1496815286624204
$ws.Range("S2").Value = 0.09235668789808917

# Row 3
$ws.Range("B3").Value = 0.006451612903225806
$ws.Range("C3").Value = 0.01290322580645161
$ws.Range("J3").Value = 0.05161290322580645
$ws.Range("P3").Value = 0.7612903225806451
$ws.Range("S3").Value = 0.167741935483871

# Row 4
$ws.Range("J4").Value = 0.07142857142857142
$ws.Range("P4").Value = 0.6714285714285714
$ws.Range("S4").Value = 0.2571428571428571

# Row 6
$ws.Range("B6").Value = 0.09090909090909091
$ws.Range("F6").Value = 0.05555555555555555
$ws.Range("J6").Value = 0.2676767676767677
$ws.Range("O6").Value = 0.0101010101010101
$ws.Range("Q6").Value = 0.1414141414141414
$ws.Range("R6").Value = 0.101010101010101
$ws.Range("S6").Value = 0.3333333333333333

# Row 7
$ws.Range("B7").Value = 0.1142857142857143
$ws.Range("D7").Value = 0.08095238095238096
$ws.Range("F7").Value = 0.06666666666666667
$ws.Range("J7").Value = 0.1523809523809524
$ws.Range("O7").Value = 0.02380952380952381
$ws.Range("Q7").Value = 0.1619047619047619
$ws.Range("R7").Value = 0.07142857142857142
$ws.Range("S7").Value = 0.3285714285714286

# Row 8
$ws.Range("B8").Value = 0.08
$ws.Range("D8").Value = 0.01894736842105263
$ws.Range("E8").Value = 0.002105263157894737
$ws.Range("F8").Value = 0.03368421052631579
$ws.Range("J8").Value = 0.1305263157894737
$ws.Range("O8").Value = 0.02105263157894737
$ws.Range("Q8").Value = 0.1873684210526316
$ws.Range("R8").Value = 0.0968421052631579
$ws.Range("S8").Value = 0.4294736842105263

# Row 9
$ws.Range("B9").Value = 0.08333333333333333
$ws.Range("D9").Value = 0.03333333333333333
$ws.Range("E9").Value = 0.005555555555555556
$ws.Range("F9").Value = 0.06111111111111111
$ws.Range("J9").Value = 0.1333333333333333
$ws.Range("O9").Value = 0.02777777777777778
$ws.Range("Q9").Value = 0.2166666666666667
$ws.Range("R9").Value = 0.08333333333333333
$ws.Range("S9").Value = 0.3555555555555556

# Row 10
$ws.Range("B10").Value = 0.1040590405904059
$ws.Range("D10").Value = 0.02730627306273063
$ws.Range("E10").Value = 0.0007380073800738007
$ws.Range("F10").Value = 0.05682656826568266
$ws.Range("J10").Value = 0.1512915129151292
$ws.Range("O10").Value = 0.00959409594095941
$ws.Range("Q10").Value = 0.2413284132841328
$ws.Range("R10").Value = 0.06863468634686347
$ws.Range("S10").Value = 0.3402214022140221

# Row 11
$ws.Range("G11").Value = 0.1335227272727273
$ws.Range("J11").Value = 0.1164772727272727
$ws.Range("K11").Value = 0.2017045454545454
$ws.Range("L11").Value = 0.5227272727272727
$ws.Range("S11").Value = 0.02556818181818182

# Row 12
$ws.Range("G12").Value = 0.6818181818181818
$ws.Range("J12").Value = 0.1919191919191919
$ws.Range("K12").Value = 0.0101010101010101
$ws.Range("L12").Value = 0.06060606060606061
$ws.Range("S12").Value = 0.05555555555555555

# Row 13
$ws.Range("G13").Value = 0.72
$ws.Range("J13").Value = 0.22
$ws.Range("S13").Value = 0.06

# Row 14
$ws.Range("G14").Value = 0.25
$ws.Range("J14").Value = 0.75

# Row 15
$ws.Range("F15").Value = 0.02192982456140351
$ws.Range("H15").Value = 0.2017543859649123
$ws.Range("I15").Value = 0.06140350877192982
$ws.Range("J15").Value = 0.3508771929824561
$ws.Range("K15").Value = 0.06140350877192982
$ws.Range("M15").Value = 0.02192982456140351
$ws.Range("O15").Value = 0.04824561403508772
$ws.Range("S15").Value = 0.2324561403508772

# Row 16
$ws.Range("F16").Value = 0.01485148514851485
$ws.Range("H16").Value = 0.1485148514851485
$ws.Range("I16").Value = 0.07920792079207921
$ws.Range("J16").Value = 0.3960396039603961
$ws.Range("K16").Value = 0.1386138613861386
$ws.Range("M16").Value = 0.02475247524752475
$ws.Range("N16").Value = 0.004950495049504951
$ws.Range("O16").Value = 0.06930693069306931
$ws.Range("S16").Value = 0.1237623762376238

# Row 17
$ws.Range("F17").Value = 0.01553398058252427
$ws.Range("H17").Value = 0.2058252427184466
$ws.Range("I17").Value = 0.05436893203883495
$ws.Range("J17").Value = 0.3766990291262136
$ws.Range("K17").Value = 0.116504854368932
$ws.Range("M17").Value = 0.01941747572815534
$ws.Range("N17").Value = 0.001941747572815534
$ws.Range("O17").Value = 0.07572815533980583
$ws.Range("S17").Value = 0.1339805825242718

# Row 18
$ws.Range("F18").Value = 0.01587301587301587
$ws.Range("H18").Value = 0.1798941798941799
$ws.Range("I18").Value = 0.05291005291005291
$ws.Range("J18").Value = 0.3333333333333333
$ws.Range("K18").Value = 0.1111111111111111
$ws.Range("M18").Value = 0.04761904761904762
$ws.Range("O18").Value = 0.09523809523809523
$ws.Range("S18").Value = 0.164021164021164

# Row 19
$ws.Range("F19").Value = 0.01632970451010887
$ws.Range("H19").Value = 0.2045101088646967
$ws.Range("I19").Value = 0.08709175738724728
$ws.Range("J19").Value = 0.3639191290824261
$ws.Range("K19").Value = 0.1135303265940902
$ws.Range("M19").Value = 0.01788491446345257
$ws.Range("N19").Value = 0.001555209953343701
$ws.Range("O19").Value = 0.06531881804043546
$ws.Range("S19").Value = 0.1298600311041991
